$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the influence experiment values in row 30 from -1 to -0.5
$ws.Range("C30").Value = -0.5
$ws.Range("D30").Value = -0.5
$ws.Range("E30").Value = -0.5
$ws.Range("F30").Value = -0.5

# Update the active selection to H32
$ws.Range("H32").Select()
